$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.472.24"
$ws.Range("E2").Value = "  +1.89%  "

# Row 3
$ws.Range("D3").Value = "1.843.69"
$ws.Range("E3").Value = "  +1.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.81%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.49%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.75%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.54%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.311"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.86%  "

# Row 11
$ws.Range("E11").Value = "  +0.87%  "

# Row 12
$ws.Range("D12").Value = "2.109.78"
$ws.Range("E12").Value = "  +1.63%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.840.87"
$ws.Range("E13").Value = "  +1.71%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.44%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.671"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.33%  "

# Row 17
$ws.Range("D17").Value = "35.443.19"
$ws.Range("E17").Value = "  +2.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.94%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  +3.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.28%  "

# Row 23
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("E24").Value = "  +3.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.122"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +33.85%  "

# Row 30
$ws.Range("E30").Value = "  +0.90%  "

# Row 31
$ws.Range("B31").Value = "EURNeutrino"
$ws.Range("C31").Value = "https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn"
$ws.Range("D31").Value = "3.352.23"
$ws.Range("E31").Value = "  +37.97%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0559"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.46%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.42%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.64%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.68%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "94.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.83%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.77%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.335.73"
$ws.Range("E39").Value = "  +1.26%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0195"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.74%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.81%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.79%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0514"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.95%  "

# Row 49
$ws.Range("D49").Value = "2.020.00"
$ws.Range("E49").Value = "  +2.30%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
